$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Activate()

# --- Update existing sample OrderId values in column R (new demo order ids) ---
$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "51490982"

$ws.Range("R3").NumberFormat = "@"
$ws.Range("R3").Value = "51490983"

$ws.Range("R5").NumberFormat = "@"
$ws.Range("R5").Value = "51490984"

$ws.Range("R6").NumberFormat = "@"
$ws.Range("R6").Value = "51490985"

# --- Add a new "PalletType" column (S) with a default "Generic Pallet" value ---
$ws.Range("S1").Value = "PalletType"
$ws.Range("S2:S13").Value = "Generic Pallet"

# Size the new column to fit its content, like the other "bestFit" columns
$ws.Columns.Item(19).AutoFit() | Out-Null

# Reflect the latest edits in the view - user left off having just retyped the OrderId column
$ws.Range("R2:R6").Select() | Out-Null
